$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "color:eye" column (D) is a categorical scale. The observation for
# accession ESP004:BGE005836 (row 2) actually carries two values separated
# by ";" (meaning it was scored as both category 1 and category 2), so we
# store that raw multi-value text in D2.
$ws.Range("D2").Value = "1;2"

# The other two accessions (rows 3 and 4) only had a single category value
# ("red" == category 1 on the categorical scale), so once the scale is
# taken into account they become plain numeric category observations.
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 1

# Row 4's "Plant size:cm" observation is refined to a decimal value.
$ws.Range("C4").Value = 3.2

# These observation cells (numeric category codes / measurements) are
# stored as text-formatted numbers, matching how the categorical/measured
# values are now represented.
$ws.Range("C2:D4").NumberFormat = "@"

# Reflect the newly reviewed range in the sheet's selection.
$ws.Range("C2:D4").Select()
